$d = $word.ActiveDocument

function Get-MatchRange([string]$text) {
    $r = $d.Content
    $ok = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $text"
    }
    return $r
}

function Split-At([int]$pos, [int]$paraEndGuess) {
    # Force a run boundary at character offset $pos by toggling a no-op
    # formatting property across [$pos, $paraEndGuess).
    $sub = $d.Range($pos, $paraEndGuess)
    $sub.Font.Bold = 1
    $sub.Font.Bold = 0
}

# ---------------------------------------------------------------------
# Paragraph 1: "Represented NJIT ... Peace Summit 2016."
# ---------------------------------------------------------------------
$r1 = Get-MatchRange("Summit 2016.")
$r1.Text = "Summit, 2016."

$p1 = Get-MatchRange("Represented NJIT as a Graduate Student Delegate at the United Nations for the Youth Leadership and Peace Summit, 2016.")
$p1Start = $p1.Start
$p1End = $p1.End

$len1 = ("Represented NJIT as a Graduate Student Delegate at the United Nations for the Yo").Length
$len2 = ("uth Leadership and Peace Summit, ").Length
$split1 = $p1Start + $len1
$split2 = $split1 + $len2

Split-At $split1 $p1End
Split-At $split2 $p1End

Write-Output "Paragraph 1 done"

# ---------------------------------------------------------------------
# Paragraph 2: "Participat" + "ed in Texas Instruments ... stage."
# The _GoBack bookmark currently sitting between the two runs here is
# relocated to paragraph 3 later, so drop it from this spot first.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$r2 = Get-MatchRange("India Design Contest 2015 and reached the quarter final stage.")
$r2.Text = "IDC and reached the quarter final stage amongst 3100 proposals,2015."

$p2 = Get-MatchRange("Participated in Texas Instruments Innovation Challenge IDC and reached the quarter final stage amongst 3100 proposals,2015.")
$p2Start = $p2.Start
$p2End = $p2.End

$runs2 = @(
    "Participat",
    "ed in Texas Instruments Innovation Challenge ID",
    "C",
    " and reached the quarter final stage",
    " amongst 3100 proposals,2015",
    "."
)
$offset = $p2Start
$splits2 = @()
for ($i = 0; $i -lt $runs2.Length - 1; $i++) {
    $offset = $offset + $runs2[$i].Length
    $splits2 += $offset
}
foreach ($sp in $splits2) {
    Split-At $sp $p2End
}

Write-Output "Paragraph 2 done"

# ---------------------------------------------------------------------
# Paragraph 3: "Head of the Association ... MIT Pune, 1st Runner's Up
# at the MIT Pune's 'Best Manager' Event."
# ---------------------------------------------------------------------
$r3a = Get-MatchRange("MIT Pune, 1")
$r3a.Text = "MIT Pune, 2015, 1"

$r3b = Get-MatchRange("Best Manager’ Event.")
$r3b.Text = "Best Manager’ Event, 2013."

$p3 = Get-MatchRange("Head of the Association of Electronics Students at MIT Pune, 2015, 1")
$p3ParaStart = $p3.Start

$p3Full = $d.Paragraphs(1)
# Re-find the whole paragraph by locating its distinctive trailing text.
$p3End = Get-MatchRange("2013.")
$paraEnd = $p3End.End

$lenA1 = ("Head of the Association of Electronics Students at MIT Pune").Length
$lenA2 = (", 2015").Length
$splitA1 = $p3ParaStart + $lenA1
$splitA2 = $splitA1 + $lenA2

Split-At $splitA1 $paraEnd
Split-At $splitA2 $paraEnd

$r3c = Get-MatchRange(" Runner’s Up at the MIT Pune’s ‘Best Manager’ Event, 2013.")
$r3cStart = $r3c.Start
$lenB1 = (" Runner’s Up at the MIT Pune’s ‘Best Manager’ Event").Length
$splitB1 = $r3cStart + $lenB1

Split-At $splitB1 $paraEnd

# Move the _GoBack bookmark here, right before the trailing period.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$finalPeriod = Get-MatchRange("2013.")
$bmPos = $finalPeriod.End - 1
$bmTarget = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmTarget)

Write-Output "Paragraph 3 done"
